$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceFormat = $ws.Range("B6").NumberFormat

# Rename existing "Schroef headers" entry to "Schroef headers 1x3" and
# flesh out its row (price, qty, supplier, URL + hyperlink).
$ws.Range("A8").Value = "Schroef headers 1x3"
$ws.Range("B8").Value = 0.59
$ws.Range("B8").NumberFormat = $priceFormat
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = "Conrad"
$ws.Range("E8").Value = "https://www.conrad.be/p/degson-dg308-254-03p-14-00ah-klemschroefblok-082-mm-aantal-polen-3-groen-1-stuks-1327224"
$ws.Hyperlinks.Add($ws.Range("E8"), "https://www.conrad.be/p/degson-dg308-254-03p-14-00ah-klemschroefblok-082-mm-aantal-polen-3-groen-1-stuks-1327224") | Out-Null

# New row: 1x2 screw headers
$ws.Range("A9").Value = "Schroef headers 1x2"
$ws.Range("B9").Value = 0.4
$ws.Range("B9").NumberFormat = $priceFormat
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "Conrad"
$ws.Range("E9").Value = "https://www.conrad.be/p/degson-dg308-254-02p-14-00ah-klemschroefblok-082-mm-aantal-polen-2-groen-1-stuks-1327242"
$ws.Hyperlinks.Add($ws.Range("E9"), "https://www.conrad.be/p/degson-dg308-254-02p-14-00ah-klemschroefblok-082-mm-aantal-polen-2-groen-1-stuks-1327242") | Out-Null

# New row: 1x6 screw headers
$ws.Range("A10").Value = "Schroef headers 1x6"
$ws.Range("B10").Value = 1.74
$ws.Range("B10").NumberFormat = $priceFormat
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "Conrad"
$ws.Range("E10").Value = "https://www.conrad.be/p/degson-dg308-254-06p-14-00ah-klemschroefblok-082-mm-aantal-polen-6-groen-1-stuks-1327217"
$ws.Hyperlinks.Add($ws.Range("E10"), "https://www.conrad.be/p/degson-dg308-254-06p-14-00ah-klemschroefblok-082-mm-aantal-polen-6-groen-1-stuks-1327217") | Out-Null

# New row: 1x8 screw headers
$ws.Range("A11").Value = "Schroef headers 1x8"
$ws.Range("B11").Value = 2.32
$ws.Range("B11").NumberFormat = $priceFormat
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Conrad"
$ws.Range("E11").Value = "https://www.conrad.be/p/degson-dg308-254-08p-14-00ah-klemschroefblok-082-mm-aantal-polen-8-groen-1-stuks-1327226"
$ws.Hyperlinks.Add($ws.Range("E11"), "https://www.conrad.be/p/degson-dg308-254-08p-14-00ah-klemschroefblok-082-mm-aantal-polen-8-groen-1-stuks-1327226") | Out-Null

$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ws.Range("E22").Select() | Out-Null
